$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.078.87'
$ws.Range('E2').Value = '  +1.34%  '
$ws.Range('D3').Value = '3.147.54'
$ws.Range('E3').Value = '  +1.99%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = "'591.09"
$ws.Range('E5').Value = '  +1.77%  '
$ws.Range('D6').Value = "'146.24"
$ws.Range('E6').Value = '  +1.28%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').Value = '3.135.04'
$ws.Range('E8').Value = '  +1.75%  '
$ws.Range('E9').Value = '  +0.78%  '
$ws.Range('E10').Value = '  +2.77%  '
$ws.Range('E11').Value = '  +5.22%  '
$ws.Range('D12').Value = "'0.457"
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('D13').Value = "'0.0000246"
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('D14').Value = "'37.14"
$ws.Range('E14').Value = '  -1.18%  '
$ws.Range('D15').Value = '3.668.25'
$ws.Range('E15').Value = '  +1.99%  '
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('E17').Value = '  +2.29%  '
$ws.Range('D18').Value = '63.865.32'
$ws.Range('E18').Value = '  +1.13%  '
$ws.Range('D19').Value = '3.141.73'
$ws.Range('E19').Value = '  +1.66%  '
$ws.Range('D20').Value = "'470.41"
$ws.Range('E20').Value = '  +2.12%  '
$ws.Range('D21').Value = "'14.39"
$ws.Range('E21').Value = '  +1.32%  '
$ws.Range('E22').Value = '  +1.15%  '
$ws.Range('D23').Value = "'7.56"
$ws.Range('E23').Value = '  +1.78%  '
$ws.Range('D24').Value = "'2.40"
$ws.Range('E24').Value = '  +13.23%  '
$ws.Range('D25').Value = "'13.19"
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('E28').Value = '  +9.94%  '
$ws.Range('E29').Value = '  +2.05%  '
$ws.Range('E30').Value = '  +7.44%  '
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('E33').Value = '  +4.34%  '
$ws.Range('D34').Value = "'27.63"
$ws.Range('E34').Value = '  +3.85%  '
$ws.Range('D35').Value = '0.0₃0857'
$ws.Range('E35').Value = '  +1.99%  '
$ws.Range('E36').Value = '  +3.26%  '
$ws.Range('E37').Value = '  +2.90%  '
$ws.Range('E38').Value = '  -0.56%  '
$ws.Range('E39').Value = '  -2.29%  '
$ws.Range('D40').Value = "'463.01"
$ws.Range('E40').Value = '  +6.91%  '
$ws.Range('D41').Value = "'51.37"
$ws.Range('E41').Value = '  +2.49%  '
$ws.Range('D42').Value = "'9.36"
$ws.Range('E42').Value = '  +7.20%  '
$ws.Range('E43').Value = '  +8.81%  '
$ws.Range('E44').Value = '  +1.44%  '
$ws.Range('D45').Value = '2.889.99'
$ws.Range('E45').Value = '  +0.97%  '
$ws.Range('D46').Value = "'40.27"
$ws.Range('E46').Value = '  +11.28%  '
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('D48').Value = "'132.55"
$ws.Range('E48').Value = '  +7.01%  '
$ws.Range('E50').Value = '  +0.90%  '
$ws.Range('E51').Value = '  +4.18%  '
